$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell so it is stored as TEXT (matching the
# source file, where these numeric-looking readings are shared-string text,
# not real numbers). Assigning a plain numeric-looking string to .Value (or
# .Formula) gets auto-coerced to a number by Excel, so instead we build a
# text-literal formula ("="value"") and convert it to a static value via
# copy / paste-special-values. This avoids minting any new NumberFormat
# style in styles.xml (unlike toggling NumberFormat to "@").
function Set-TextValue($range, [string]$text) {
    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# Copy the existing date cell's format (style) down into the two new rows so
# we reuse the already-defined date style (s="1") instead of minting a new
# one.
$ws.Range("A100").Copy($ws.Range("A101:A102"))

# --- Row 101 (2025-04-01) ---
$ws.Range("A101").Value = 45748
$ws.Range("B101").Value = 133.091761295188
$ws.Range("C101").Value = 123.545536029778
Set-TextValue $ws.Range("D101") "113.3"
Set-TextValue $ws.Range("E101") "115.0"
Set-TextValue $ws.Range("F101") " 88.1"
Set-TextValue $ws.Range("G101") "172.2"

# --- Row 102 (2025-05-01) ---
$ws.Range("A102").Value = 45778
$ws.Range("B102").Value = 132.510366166753
$ws.Range("C102").Value = 124.395161380377
Set-TextValue $ws.Range("D102") "114.9"
Set-TextValue $ws.Range("E102") "116.0"
Set-TextValue $ws.Range("F102") " 88.8"
Set-TextValue $ws.Range("G102") "172.8"
